$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking price strings
# (e.g. "1.00", "0.630") are preserved exactly as text and not coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.978.09"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "3.564.81"

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "592.54"
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").Value = "197.49"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  -2.15%  "

$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("D10").Value = "0.630"
$ws.Range("E10").Value = "  -2.75%  "

$ws.Range("D11").Value = "53.24"
$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").Value = "0.0000291"
$ws.Range("E12").Value = "  -4.10%  "

$ws.Range("D13").Value = "9.33"
$ws.Range("E13").Value = "  -2.50%  "

$ws.Range("D14").Value = "4.136.35"
$ws.Range("E14").Value = "  -1.44%  "

$ws.Range("D15").Value = "654.20"
$ws.Range("E15").Value = "  +8.56%  "

$ws.Range("D16").Value = "70.014.08"
$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").Value = "12.74"
$ws.Range("E17").Value = "  -2.09%  "

$ws.Range("D18").Value = "3.572.12"
$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "18.53"
$ws.Range("E19").Value = "  -2.73%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "0.122"
$ws.Range("E20").Value = "  -0.97%  "

$ws.Range("D21").Value = "0.969"
$ws.Range("E21").Value = "  -2.81%  "

$ws.Range("D22").Value = "18.39"
$ws.Range("E22").Value = "  +3.05%  "

$ws.Range("D23").Value = "5.43"
$ws.Range("E23").Value = "  +4.57%  "

$ws.Range("D24").Value = "104.82"
$ws.Range("E24").Value = "  +2.99%  "

$ws.Range("D25").Value = "4.43"
$ws.Range("E25").Value = "  -4.26%  "

$ws.Range("E26").Value = "  -2.52%  "

$ws.Range("D27").Value = "10.32"
$ws.Range("E27").Value = "  -4.13%  "

$ws.Range("D28").Value = "9.68"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").Value = "33.45"
$ws.Range("E29").Value = "  -1.20%  "

$ws.Range("D30").Value = "4.36"

$ws.Range("D31").Value = "6.86"
$ws.Range("E31").Value = "  -5.22%  "

$ws.Range("D32").Value = "11.84"
$ws.Range("E32").Value = "  -3.66%  "

$ws.Range("E33").Value = "  -5.24%  "

$ws.Range("D34").Value = "61.99"
$ws.Range("E34").Value = "  -2.19%  "

$ws.Range("D35").Value = "3.757.00"
$ws.Range("E35").Value = "  -3.62%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0819"
$ws.Range("E36").Value = "  -7.52%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "3.69"
$ws.Range("E38").Value = "  +3.85%  "

$ws.Range("D39").Value = "514.75"
$ws.Range("E39").Value = "  -5.01%  "

$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  -4.09%  "

$ws.Range("D41").Value = "0.375"
$ws.Range("E41").Value = "  -3.98%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "35.21"
$ws.Range("E42").Value = "  -4.80%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.135"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("D44").Value = "0.0454"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("D45").Value = "3.42"
$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  +0.95%  "

$ws.Range("D47").Value = "0.138"
$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("E48").Value = "  -0.19%  "

$ws.Range("D49").Value = "8.39"
$ws.Range("E49").Value = "  -2.51%  "

$ws.Range("E50").Value = "  +18.24%  "

$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "0.000240"
$ws.Range("E51").Value = "  -4.59%  "

# Restore the default (unstyled) look for column D now that the text is set,
# matching the workbook's original cell styling (no explicit style index).
$ws.Range("D2:D51").Style = "Normal"
